$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3128.5
$ws.Range("I113").Value = 2900
$ws.Range("J113").Value = 3219.9
$ws.Range("K113").Value = 2900
$ws.Range("L113").Value = 3219.9
$ws.Range("M113").Value = 354
$ws.Range("N113").Value = -9727.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 878.8823
$ws.Range("I2").Value = 871.5
$ws.Range("J2").Value = 913.3333
$ws.Range("K2").Value = 871.5
$ws.Range("L2").Value = 913.3333
$ws.Range("M2").Value = -758.5
$ws.Range("N2").Value = -1139.3333
$ws.Range("H32").Value = 16288.755
$ws.Range("I32").Value = 4291.7144
$ws.Range("J32").Value = 37883.43
$ws.Range("K32").Value = 4291.7144
$ws.Range("L32").Value = 37883.43
$ws.Range("M32").Value = -4004.7144
$ws.Range("N32").Value = -38457.43
$ws.Range("H45").Value = 2430.1667
$ws.Range("J45").Value = 1960
$ws.Range("L45").Value = 1960
$ws.Range("N45").Value = -2714
$ws.Range("H63").Value = 2847.5264
$ws.Range("I63").Value = 2147.8572
$ws.Range("J63").Value = 4806.6
$ws.Range("K63").Value = 2147.8572
$ws.Range("L63").Value = 4806.6
$ws.Range("M63").Value = -1461.8572
$ws.Range("N63").Value = -6178.6
$ws.Range("H66").Value = 2847.5264
$ws.Range("I66").Value = 2147.8572
$ws.Range("J66").Value = 4806.6
$ws.Range("K66").Value = 10739.286
$ws.Range("L66").Value = 24033
$ws.Range("M66").Value = -7307.286
$ws.Range("N66").Value = -30897
$ws.Range("H116").Value = 878.8823
$ws.Range("I116").Value = 871.5
$ws.Range("J116").Value = 913.3333
$ws.Range("K116").Value = 871.5
$ws.Range("L116").Value = 913.3333
$ws.Range("M116").Value = 1422.5
$ws.Range("N116").Value = -5501.3333
$ws.Range("H132").Value = 1633.5
$ws.Range("I132").Value = 877.89655
$ws.Range("J132").Value = 2676.9524
$ws.Range("K132").Value = 2633.68965
$ws.Range("L132").Value = 8030.8572
$ws.Range("M132").Value = -103.6896500000003
$ws.Range("N132").Value = -13090.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 878.8823
$ws.Range("I3").Value = 871.5
$ws.Range("J3").Value = 913.3333
$ws.Range("K3").Value = 871.5
$ws.Range("L3").Value = 913.3333
$ws.Range("M3").Value = -757.5
$ws.Range("N3").Value = -1141.3333
$ws.Range("H86").Value = 28572548
$ws.Range("I86").Value = 28572548
$ws.Range("K86").Value = 28572548
$ws.Range("M86").Value = -28571425
$ws.Range("H89").Value = 28572548
$ws.Range("I89").Value = 28572548
$ws.Range("K89").Value = 142862740
$ws.Range("M89").Value = -142857124

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 14727.777
$ws.Range("J80").Value = 14727.777
$ws.Range("L80").Value = 14727.777
$ws.Range("N80").Value = -16973.777
$ws.Range("H83").Value = 14727.777
$ws.Range("J83").Value = 14727.777
$ws.Range("L83").Value = 44183.331
$ws.Range("N83").Value = -55415.331
$ws.Range("H99").Value = 3796.8262
$ws.Range("I99").Value = 3945.4375
$ws.Range("J99").Value = 3457.1428
$ws.Range("K99").Value = 3945.4375
$ws.Range("L99").Value = 3457.1428
$ws.Range("M99").Value = -2447.4375
$ws.Range("N99").Value = -6453.1428
$ws.Range("H107").Value = 6400
$ws.Range("I107").Value = 9800
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 9800
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -7880
$ws.Range("N107").Value = -6840
$ws.Range("H126").Value = 3796.8262
$ws.Range("I126").Value = 3945.4375
$ws.Range("J126").Value = 3457.1428
$ws.Range("K126").Value = 11836.3125
$ws.Range("L126").Value = 10371.4284
$ws.Range("M126").Value = -9366.3125
$ws.Range("N126").Value = -15311.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 810.3333
$ws.Range("I18").Value = 810.3333
$ws.Range("K18").Value = 2430.9999
$ws.Range("M18").Value = -2261.9999
$ws.Range("H131").Value = 1905845.4
$ws.Range("I131").Value = 11111640
$ws.Range("J131").Value = 1198.1034
$ws.Range("K131").Value = 33334920
$ws.Range("L131").Value = 3594.3102
$ws.Range("M131").Value = -33329880
$ws.Range("N131").Value = -13674.3102
$ws.Range("H133").Value = 1166.6666
$ws.Range("J133").Value = 1300
$ws.Range("L133").Value = 3900
$ws.Range("N133").Value = -14020
$ws.Range("H134").Value = 1895.2174
$ws.Range("I134").Value = 3011.25
$ws.Range("K134").Value = 9033.75
$ws.Range("M134").Value = -3963.75
$ws.Range("H137").Value = 5507.3706
$ws.Range("I137").Value = 784
$ws.Range("J137").Value = 7160.55
$ws.Range("K137").Value = 2352
$ws.Range("L137").Value = 21481.65
$ws.Range("M137").Value = 2748
$ws.Range("N137").Value = -31681.65
$ws.Range("H138").Value = 1084.5454
$ws.Range("I138").Value = 1915
$ws.Range("K138").Value = 5745
$ws.Range("M138").Value = -605
$ws.Range("H139").Value = 4684.8945
$ws.Range("I139").Value = 2987.5715
$ws.Range("J139").Value = 5675
$ws.Range("K139").Value = 8962.7145
$ws.Range("L139").Value = 17025
$ws.Range("M139").Value = -3822.7145
$ws.Range("N139").Value = -27305

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2467.2727
$ws.Range("I80").Value = 2471.4285
$ws.Range("J80").Value = 2465.3333
$ws.Range("K80").Value = 2471.4285
$ws.Range("L80").Value = 2465.3333
$ws.Range("M80").Value = -1473.4285
$ws.Range("N80").Value = -4461.3333
$ws.Range("H83").Value = 2467.2727
$ws.Range("I83").Value = 2471.4285
$ws.Range("J83").Value = 2465.3333
$ws.Range("K83").Value = 12357.1425
$ws.Range("L83").Value = 12326.6665
$ws.Range("M83").Value = -7365.1425
$ws.Range("N83").Value = -22310.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2376.1177
$ws.Range("I7").Value = 2254.1428
$ws.Range("K7").Value = 2254.1428
$ws.Range("M7").Value = -2142.1428
$ws.Range("H82").Value = 1387.5333
$ws.Range("I82").Value = 776.625
$ws.Range("J82").Value = 2085.7144
$ws.Range("K82").Value = 776.625
$ws.Range("L82").Value = 2085.7144
$ws.Range("M82").Value = -415.625
$ws.Range("N82").Value = -2807.7144
$ws.Range("H85").Value = 1387.5333
$ws.Range("I85").Value = 776.625
$ws.Range("J85").Value = 2085.7144
$ws.Range("K85").Value = 776.625
$ws.Range("L85").Value = 2085.7144
$ws.Range("M85").Value = 471.375
$ws.Range("N85").Value = -4581.7144
$ws.Range("H126").Value = 2376.1177
$ws.Range("I126").Value = 2254.1428
$ws.Range("K126").Value = 6762.428400000001
$ws.Range("M126").Value = -4292.428400000001
$ws.Range("H136").Value = 4677.8647
$ws.Range("I136").Value = 7972.3335
$ws.Range("J136").Value = 1556.7894
$ws.Range("K136").Value = 23917.0005
$ws.Range("L136").Value = 4670.3682
$ws.Range("M136").Value = -21367.0005
$ws.Range("N136").Value = -9770.368200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 55557856
$ws.Range("I81").Value = 90911224
$ws.Range("K81").Value = 181822448
$ws.Range("M81").Value = -181821387
$ws.Range("H84").Value = 55557856
$ws.Range("I84").Value = 90911224
$ws.Range("K84").Value = 909112240
$ws.Range("M84").Value = -909106936
$ws.Range("H126").Value = 2576.1428
$ws.Range("I126").Value = 4686
$ws.Range("J126").Value = 1277.7693
$ws.Range("K126").Value = 14058
$ws.Range("L126").Value = 3833.3079
$ws.Range("M126").Value = -11588
$ws.Range("N126").Value = -8773.3079
$ws.Range("H136").Value = 777.1429000000001
$ws.Range("I136").Value = 490.69766
$ws.Range("J136").Value = 2830
$ws.Range("K136").Value = 1472.09298
$ws.Range("L136").Value = 8490
$ws.Range("M136").Value = 1077.90702
$ws.Range("N136").Value = -13590
